# Updated cryptos list on Mon Sep 16 11:32:26 UTC 2024 with GitHub Actions
#
# The sheet stores "Price" (col D) and "Volume(1h)" (col E) as plain text
# (not numbers) so values like "58.784.49" (dotted thousands) and
# "  -1.94%  " (padded percent strings) round-trip byte-for-byte.
#
# Writing a numeric-looking literal straight into Range.Value would make
# Excel auto-detect it as a Number and silently reformat it (e.g.
# "0.0500" -> 0.05). To keep those cells as Text - matching the source
# workbook - a leading apostrophe (the same trick used in the Excel UI)
# forces text entry for values that would otherwise parse as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    # Plain decimal number (optionally signed) -> Excel would auto-convert
    # Range.Value to a Number and reformat it. Dotted-thousands strings like
    # "58.720.62" or padded percents like "  -2.07%  " never match this and
    # are always safe to assign as-is.
    if ($Text.Trim() -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Leading apostrophe -> force Text type, exactly like typing it in Excel.
        $ws.Range($Cell).Value = "'" + $Text
    } else {
        $ws.Range($Cell).Value = $Text
    }
}

# r,c -> cell ref : (Price/D, Volume(1h)/E) new values, per the diff.
Set-TextValue "D2"  "58.720.62"
Set-TextValue "E2"  "  -2.07%  "

Set-TextValue "D3"  "2.303.01"

Set-TextValue "E4"  "  +0.03%  "

Set-TextValue "D5"  "547.05"
Set-TextValue "E5"  "  -1.19%  "

Set-TextValue "D6"  "131.93"
Set-TextValue "E6"  "  -2.78%  "

Set-TextValue "E7"  "  +0.06%  "

Set-TextValue "D8"  "0.571"
Set-TextValue "E8"  "  -2.06%  "

Set-TextValue "D9"  "2.302.67"
Set-TextValue "E9"  "  -4.25%  "

Set-TextValue "E10" "  -3.09%  "

Set-TextValue "D11" "5.51"
Set-TextValue "E11" "  -1.94%  "

Set-TextValue "E12" "  +1.36%  "

Set-TextValue "E13" "  -4.67%  "

Set-TextValue "E14" "  -3.03%  "

Set-TextValue "D15" "2.713.85"
Set-TextValue "E15" "  -4.37%  "

Set-TextValue "D16" "58.698.06"
Set-TextValue "E16" "  -1.90%  "

Set-TextValue "E17" "  -3.00%  "

Set-TextValue "D18" "2.341.68"
Set-TextValue "E18" "  -0.05%  "

Set-TextValue "D19" "10.68"
Set-TextValue "E19" "  -4.35%  "

Set-TextValue "E20" "  -4.12%  "

Set-TextValue "D21" "314.56"
Set-TextValue "E21" "  -3.61%  "

Set-TextValue "D22" "6.47"
Set-TextValue "E22" "  -4.20%  "

Set-TextValue "E23" "  -0.02%  "

Set-TextValue "D24" "63.42"
Set-TextValue "E24" "  -1.85%  "

Set-TextValue "D25" "0.168"

Set-TextValue "E26" "  +0.01%  "

Set-TextValue "E27" "  -5.61%  "

Set-TextValue "D28" "1.32"
Set-TextValue "E28" "  -5.46%  "

Set-TextValue "D29" "1.76"
Set-TextValue "E29" "  -1.73%  "

Set-TextValue "D30" "168.51"
Set-TextValue "E30" "  -0.84%  "

Set-TextValue "E31" "  -5.13%  "

Set-TextValue "E32" "  +0.78%  "

Set-TextValue "E33" "  -5.26%  "

Set-TextValue "E34" "  -4.60%  "

Set-TextValue "E35" "  -0.01%  "

Set-TextValue "E36" "  -3.24%  "

Set-TextValue "E37" "  -0.02%  "

Set-TextValue "E38" "  -4.31%  "

Set-TextValue "D39" "3.97"
Set-TextValue "E39" "  -5.02%  "

Set-TextValue "D40" "38.07"
Set-TextValue "E40" "  -1.21%  "

Set-TextValue "E41" "  -4.81%  "

Set-TextValue "D42" "297.88"
Set-TextValue "E42" "  -7.62%  "

Set-TextValue "D43" "141.26"
Set-TextValue "E43" "  -3.52%  "

Set-TextValue "E44" "  -4.16%  "

Set-TextValue "D45" "0.0951"
Set-TextValue "E45" "  -1.13%  "

Set-TextValue "D46" "0.0500"
Set-TextValue "E46" "  -2.48%  "

Set-TextValue "D47" "0.555"
Set-TextValue "E47" "  -3.30%  "

Set-TextValue "D48" "18.47"
Set-TextValue "E48" "  -6.85%  "

Set-TextValue "E49" "  -2.62%  "

Set-TextValue "D50" "16.62"
Set-TextValue "E50" "  -3.63%  "

Set-TextValue "D51" "11.02"
Set-TextValue "E51" "  -0.20%  "
